{"js": "const pairs = [\n  [\"N = 16,490\", \"N = 3,264\"],\n  [\"N = 9,362\", \"N = 1,911\"],\n  [\"58.2 \u00b1 18.1; 0.0/16,490.0 missing (0.0%)\", \"58.4 \u00b1 17.7; 0.0/3,264.0 missing (0.0%)\"],\n  [\"61.1 \u00b1 16.9; 0.0/9,362.0 missing (0.0%)\", \"61.3 \u00b1 17.2; 0.0/1,911.0 missing (0.0%)\"],\n  [\"32.5 \u00b1 8.7; 9,223.0/16,490.0 missing (55.9%)\", \"32.6 \u00b1 8.9; 1,821.0/3,264.0 missing (55.8%)\"],\n  [\"29.0 \u00b1 7.1; 5,500.0/9,362.0 missing (58.7%)\", \"28.8 \u00b1 7.0; 1,113.0/1,911.0 missing (58.2%)\"],\n  [\"8,598 (52%)\", \"1,737 (53%)\"],\n  [\"4,218 (45%)\", \"854 (45%)\"],\n  [\"7,892 (48%)\", \"1,527 (47%)\"],\n  [\"5,144 (55%)\", \"1,057 (55%)\"],\n  [\"10,061 (61%)\", \"2,009 (62%)\"],\n  [\"6,090 (65%)\", \"1,271 (67%)\"],\n  [\"3,149 (19%)\", \"601 (18%)\"],\n  [\"1,387 (15%)\", \"283 (15%)\"],\n  [\"1,168 (7.1%)\", \"250 (7.7%)\"],\n  [\"500 (5.3%)\", \"84 (4.4%)\"],\n  [\"252 (1.5%)\", \"46 (1.4%)\"],\n  [\"179 (1.9%)\", \"37 (1.9%)\"],\n  [\"88 (0.5%)\", \"11 (0.3%)\"],\n  [\"101 (1.1%)\", \"25 (1.3%)\"],\n  [\"25 (0.2%)\", \"8 (0.2%)\"],\n  [\"14 (0.1%)\", \"2 (0.1%)\"],\n  [\"1,747 (11%)\", \"339 (10%)\"],\n  [\"1,091 (12%)\", \"209 (11%)\"],\n  [\"6,962 (42%)\", \"1,374 (42%)\"],\n  [\"5,171 (55%)\", \"1,057 (55%)\"],\n  [\"4,707 (29%)\", \"921 (28%)\"],\n  [\"1,831 (20%)\", \"379 (20%)\"],\n  [\"1,143 (6.9%)\", \"246 (7.5%)\"],\n  [\"783 (8.4%)\", \"160 (8.4%)\"],\n  [\"3,678 (22%)\", \"723 (22%)\"],\n  [\"1,577 (17%)\", \"315 (16%)\"],\n  [\"3,047 (18%)\", \"576 (18%)\"],\n  [\"1,455 (16%)\", \"303 (16%)\"],\n  [\"2,439 (15%)\", \"464 (14%)\"],\n  [\"1,086 (12%)\", \"220 (12%)\"],\n  [\"3,037 (18%)\", \"579 (18%)\"],\n  [\"2,052 (22%)\", \"434 (23%)\"],\n  [\"2,964 (18%)\", \"587 (18%)\"],\n  [\"2,104 (22%)\", \"457 (24%)\"],\n  [\"599 (3.6%)\", \"130 (4.0%)\"],\n  [\"428 (4.6%)\", \"103 (5.4%)\"],\n  [\"1,231 (7.5%)\", \"239 (7.3%)\"],\n  [\"897 (9.6%)\", \"195 (10%)\"],\n  [\"2,825 (17%)\", \"519 (16%)\"],\n  [\"1,837 (20%)\", \"385 (20%)\"],\n  [\"4,838 (29%)\", \"954 (29%)\"],\n  [\"2,776 (30%)\", \"559 (29%)\"],\n  [\"7,055 (43%)\", \"1,427 (44%)\"],\n  [\"1,451 (15%)\", \"283 (15%)\"],\n  [\"9,435 (57%)\", \"1,837 (56%)\"],\n  [\"7,911 (85%)\", \"1,628 (85%)\"],\n  [\"N = 18,392\", \"N = 3,730\"],\n  [\"N = 7,460\", \"N = 1,445\"],\n  [\"59.4 \u00b1 17.7; 0.0/18,392.0 missing (0.0%)\", \"59.6 \u00b1 17.6; 0.0/3,730.0 missing (0.0%)\"],\n  [\"58.9 \u00b1 17.7; 0.0/7,460.0 missing (0.0%)\", \"59.0 \u00b1 17.7; 0.0/1,445.0 missing (0.0%)\"],\n  [\"31.9 \u00b1 8.5; 9,691.0/18,392.0 missing (52.7%)\", \"31.9 \u00b1 8.6; 1,939.0/3,730.0 missing (52.0%)\"],\n  [\"29.0 \u00b1 7.5; 5,032.0/7,460.0 missing (67.5%)\", \"28.9 \u00b1 7.5; 995.0/1,445.0 missing (68.9%)\"],\n  [\"9,331 (51%)\", \"1,921 (52%)\"],\n  [\"3,485 (47%)\", \"670 (46%)\"],\n  [\"9,061 (49%)\", \"1,809 (48%)\"],\n  [\"3,975 (53%)\", \"775 (54%)\"],\n  [\"12,190 (66%)\", \"2,513 (67%)\"],\n  [\"3,961 (53%)\", \"767 (53%)\"],\n  [\"3,109 (17%)\", \"616 (17%)\"],\n  [\"1,427 (19%)\", \"268 (19%)\"],\n  [\"1,149 (6.2%)\", \"221 (5.9%)\"],\n  [\"519 (7.0%)\", \"113 (7.8%)\"],\n  [\"279 (1.5%)\", \"54 (1.4%)\"],\n  [\"152 (2.0%)\", \"29 (2.0%)\"],\n  [\"91 (0.5%)\", \"13 (0.3%)\"],\n  [\"98 (1.3%)\", \"23 (1.6%)\"],\n  [\"33 (0.2%)\", \"8 (0.2%)\"],\n  [\"6 (<0.1%)\", \"2 (0.1%)\"],\n  [\"1,541 (8.4%)\", \"305 (8.2%)\"],\n  [\"1,297 (17%)\", \"243 (17%)\"],\n  [\"9,882 (54%)\", \"1,987 (53%)\"],\n  [\"2,251 (30%)\", \"444 (31%)\"],\n  [\"3,358 (18%)\", \"703 (19%)\"],\n  [\"3,180 (43%)\", \"597 (41%)\"],\n  [\"1,241 (6.7%)\", \"264 (7.1%)\"],\n  [\"685 (9.2%)\", \"142 (9.8%)\"],\n  [\"3,911 (21%)\", \"776 (21%)\"],\n  [\"1,344 (18%)\", \"262 (18%)\"],\n  [\"3,273 (18%)\", \"648 (17%)\"],\n  [\"1,229 (16%)\", \"231 (16%)\"],\n  [\"2,520 (14%)\", \"513 (14%)\"],\n  [\"1,005 (13%)\", \"171 (12%)\"],\n  [\"3,596 (20%)\", \"731 (20%)\"],\n  [\"1,493 (20%)\", \"282 (20%)\"],\n  [\"3,495 (19%)\", \"741 (20%)\"],\n  [\"1,573 (21%)\", \"303 (21%)\"],\n  [\"757 (4.1%)\", \"175 (4.7%)\"],\n  [\"270 (3.6%)\", \"58 (4.0%)\"],\n  [\"1,405 (7.6%)\", \"295 (7.9%)\"],\n  [\"723 (9.7%)\", \"139 (9.6%)\"],\n  [\"3,119 (17%)\", \"628 (17%)\"],\n  [\"1,543 (21%)\", \"276 (19%)\"],\n  [\"5,137 (28%)\", \"1,034 (28%)\"],\n  [\"2,477 (33%)\", \"479 (33%)\"],\n  [\"6,183 (34%)\", \"1,247 (33%)\"],\n  [\"2,323 (31%)\", \"463 (32%)\"],\n  [\"12,209 (66%)\", \"2,483 (67%)\"],\n  [\"5,137 (69%)\", \"982 (68%)\"],\n];\n\nconst body = context.document.body;\nlet totalReplaced = 0;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n    totalReplaced++;\n  }\n  await context.sync();\n}\nreturn totalReplaced;", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"N = 16,490\", \"N = 3,264\")\n  ,@(\"N = 9,362\", \"N = 1,911\")\n  ,@(\"58.2 \u00b1 18.1; 0.0/16,490.0 missing (0.0%)\", \"58.4 \u00b1 17.7; 0.0/3,264.0 missing (0.0%)\")\n  ,@(\"61.1 \u00b1 16.9; 0.0/9,362.0 missing (0.0%)\", \"61.3 \u00b1 17.2; 0.0/1,911.0 missing (0.0%)\")\n  ,@(\"32.5 \u00b1 8.7; 9,223.0/16,490.0 missing (55.9%)\", \"32.6 \u00b1 8.9; 1,821.0/3,264.0 missing (55.8%)\")\n  ,@(\"29.0 \u00b1 7.1; 5,500.0/9,362.0 missing (58.7%)\", \"28.8 \u00b1 7.0; 1,113.0/1,911.0 missing (58.2%)\")\n  ,@(\"8,598 (52%)\", \"1,737 (53%)\")\n  ,@(\"4,218 (45%)\", \"854 (45%)\")\n  ,@(\"7,892 (48%)\", \"1,527 (47%)\")\n  ,@(\"5,144 (55%)\", \"1,057 (55%)\")\n  ,@(\"10,061 (61%)\", \"2,009 (62%)\")\n  ,@(\"6,090 (65%)\", \"1,271 (67%)\")\n  ,@(\"3,149 (19%)\", \"601 (18%)\")\n  ,@(\"1,387 (15%)\", \"283 (15%)\")\n  ,@(\"1,168 (7.1%)\", \"250 (7.7%)\")\n  ,@(\"500 (5.3%)\", \"84 (4.4%)\")\n  ,@(\"252 (1.5%)\", \"46 (1.4%)\")\n  ,@(\"179 (1.9%)\", \"37 (1.9%)\")\n  ,@(\"88 (0.5%)\", \"11 (0.3%)\")\n  ,@(\"101 (1.1%)\", \"25 (1.3%)\")\n  ,@(\"25 (0.2%)\", \"8 (0.2%)\")\n  ,@(\"14 (0.1%)\", \"2 (0.1%)\")\n  ,@(\"1,747 (11%)\", \"339 (10%)\")\n  ,@(\"1,091 (12%)\", \"209 (11%)\")\n  ,@(\"6,962 (42%)\", \"1,374 (42%)\")\n  ,@(\"5,171 (55%)\", \"1,057 (55%)\")\n  ,@(\"4,707 (29%)\", \"921 (28%)\")\n  ,@(\"1,831 (20%)\", \"379 (20%)\")\n  ,@(\"1,143 (6.9%)\", \"246 (7.5%)\")\n  ,@(\"783 (8.4%)\", \"160 (8.4%)\")\n  ,@(\"3,678 (22%)\", \"723 (22%)\")\n  ,@(\"1,577 (17%)\", \"315 (16%)\")\n  ,@(\"3,047 (18%)\", \"576 (18%)\")\n  ,@(\"1,455 (16%)\", \"303 (16%)\")\n  ,@(\"2,439 (15%)\", \"464 (14%)\")\n  ,@(\"1,086 (12%)\", \"220 (12%)\")\n  ,@(\"3,037 (18%)\", \"579 (18%)\")\n  ,@(\"2,052 (22%)\", \"434 (23%)\")\n  ,@(\"2,964 (18%)\", \"587 (18%)\")\n  ,@(\"2,104 (22%)\", \"457 (24%)\")\n  ,@(\"599 (3.6%)\", \"130 (4.0%)\")\n  ,@(\"428 (4.6%)\", \"103 (5.4%)\")\n  ,@(\"1,231 (7.5%)\", \"239 (7.3%)\")\n  ,@(\"897 (9.6%)\", \"195 (10%)\")\n  ,@(\"2,825 (17%)\", \"519 (16%)\")\n  ,@(\"1,837 (20%)\", \"385 (20%)\")\n  ,@(\"4,838 (29%)\", \"954 (29%)\")\n  ,@(\"2,776 (30%)\", \"559 (29%)\")\n  ,@(\"7,055 (43%)\", \"1,427 (44%)\")\n  ,@(\"1,451 (15%)\", \"283 (15%)\")\n  ,@(\"9,435 (57%)\", \"1,837 (56%)\")\n  ,@(\"7,911 (85%)\", \"1,628 (85%)\")\n  ,@(\"N = 18,392\", \"N = 3,730\")\n  ,@(\"N = 7,460\", \"N = 1,445\")\n  ,@(\"59.4 \u00b1 17.7; 0.0/18,392.0 missing (0.0%)\", \"59.6 \u00b1 17.6; 0.0/3,730.0 missing (0.0%)\")\n  ,@(\"58.9 \u00b1 17.7; 0.0/7,460.0 missing (0.0%)\", \"59.0 \u00b1 17.7; 0.0/1,445.0 missing (0.0%)\")\n  ,@(\"31.9 \u00b1 8.5; 9,691.0/18,392.0 missing (52.7%)\", \"31.9 \u00b1 8.6; 1,939.0/3,730.0 missing (52.0%)\")\n  ,@(\"29.0 \u00b1 7.5; 5,032.0/7,460.0 missing (67.5%)\", \"28.9 \u00b1 7.5; 995.0/1,445.0 missing (68.9%)\")\n  ,@(\"9,331 (51%)\", \"1,921 (52%)\")\n  ,@(\"3,485 (47%)\", \"670 (46%)\")\n  ,@(\"9,061 (49%)\", \"1,809 (48%)\")\n  ,@(\"3,975 (53%)\", \"775 (54%)\")\n  ,@(\"12,190 (66%)\", \"2,513 (67%)\")\n  ,@(\"3,961 (53%)\", \"767 (53%)\")\n  ,@(\"3,109 (17%)\", \"616 (17%)\")\n  ,@(\"1,427 (19%)\", \"268 (19%)\")\n  ,@(\"1,149 (6.2%)\", \"221 (5.9%)\")\n  ,@(\"519 (7.0%)\", \"113 (7.8%)\")\n  ,@(\"279 (1.5%)\", \"54 (1.4%)\")\n  ,@(\"152 (2.0%)\", \"29 (2.0%)\")\n  ,@(\"91 (0.5%)\", \"13 (0.3%)\")\n  ,@(\"98 (1.3%)\", \"23 (1.6%)\")\n  ,@(\"33 (0.2%)\", \"8 (0.2%)\")\n  ,@(\"6 (<0.1%)\", \"2 (0.1%)\")\n  ,@(\"1,541 (8.4%)\", \"305 (8.2%)\")\n  ,@(\"1,297 (17%)\", \"243 (17%)\")\n  ,@(\"9,882 (54%)\", \"1,987 (53%)\")\n  ,@(\"2,251 (30%)\", \"444 (31%)\")\n  ,@(\"3,358 (18%)\", \"703 (19%)\")\n  ,@(\"3,180 (43%)\", \"597 (41%)\")\n  ,@(\"1,241 (6.7%)\", \"264 (7.1%)\")\n  ,@(\"685 (9.2%)\", \"142 (9.8%)\")\n  ,@(\"3,911 (21%)\", \"776 (21%)\")\n  ,@(\"1,344 (18%)\", \"262 (18%)\")\n  ,@(\"3,273 (18%)\", \"648 (17%)\")\n  ,@(\"1,229 (16%)\", \"231 (16%)\")\n  ,@(\"2,520 (14%)\", \"513 (14%)\")\n  ,@(\"1,005 (13%)\", \"171 (12%)\")\n  ,@(\"3,596 (20%)\", \"731 (20%)\")\n  ,@(\"1,493 (20%)\", \"282 (20%)\")\n  ,@(\"3,495 (19%)\", \"741 (20%)\")\n  ,@(\"1,573 (21%)\", \"303 (21%)\")\n  ,@(\"757 (4.1%)\", \"175 (4.7%)\")\n  ,@(\"270 (3.6%)\", \"58 (4.0%)\")\n  ,@(\"1,405 (7.6%)\", \"295 (7.9%)\")\n  ,@(\"723 (9.7%)\", \"139 (9.6%)\")\n  ,@(\"3,119 (17%)\", \"628 (17%)\")\n  ,@(\"1,543 (21%)\", \"276 (19%)\")\n  ,@(\"5,137 (28%)\", \"1,034 (28%)\")\n  ,@(\"2,477 (33%)\", \"479 (33%)\")\n  ,@(\"6,183 (34%)\", \"1,247 (33%)\")\n  ,@(\"2,323 (31%)\", \"463 (32%)\")\n  ,@(\"12,209 (66%)\", \"2,483 (67%)\")\n  ,@(\"5,137 (69%)\", \"982 (68%)\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $ok = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $ok) {\n    throw \"Not found: $oldText\"\n  }\n}"}
